# Actualiza el diccionario de datos (data dictionary) on sheet "Hoja1".
#
#   M1: "Fecha_obligación"  -> "Fecha_obligacion"            (accent removed)
#   S1: "NIT_Beneficiario"  -> "Identificacion_Beneficiario" (field renamed)
#
# The data row (row 2) is untouched - only the header labels change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M1").Value = "Fecha_obligacion"
$ws.Range("S1").Value = "Identificacion_Beneficiario"

# Cosmetic view state captured alongside the dictionary edit.
$ws.Application.ActiveWindow.Zoom = 115
$ws.Range("S4").Select()
